$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 78100.766
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 78100.766
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 234302.298
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -234638.298

$ws.Range("H40").Value = 1538.4615
$ws.Range("I40").Value = 1425
$ws.Range("J40").Value = 1720
$ws.Range("K40").Value = 1425
$ws.Range("L40").Value = 1720
$ws.Range("M40").Value = -1250
$ws.Range("N40").Value = -2070

$ws.Range("H58").Value = 2201.913
$ws.Range("I58").Value = 349
$ws.Range("J58").Value = 2716.611
$ws.Range("K58").Value = 1047
$ws.Range("L58").Value = 8149.833
$ws.Range("M58").Value = -897
$ws.Range("N58").Value = -8449.832999999999

$ws.Range("H92").Value = 981.5789
$ws.Range("I92").Value = 1073.125
$ws.Range("J92").Value = 493.33334
$ws.Range("K92").Value = 1073.125
$ws.Range("L92").Value = 493.33334
$ws.Range("M92").Value = 174.875
$ws.Range("N92").Value = -2989.33334

$ws.Range("H98").Value = 8666.333000000001
$ws.Range("I98").Value = 5583.75
$ws.Range("J98").Value = 20996.666
$ws.Range("K98").Value = 5583.75
$ws.Range("L98").Value = 20996.666
$ws.Range("M98").Value = -4085.75
$ws.Range("N98").Value = -23992.666

$ws.Range("H122").Value = 8666.333000000001
$ws.Range("I122").Value = 5583.75
$ws.Range("J122").Value = 20996.666
$ws.Range("K122").Value = 16751.25
$ws.Range("L122").Value = 62989.99800000001
$ws.Range("M122").Value = -14301.25
$ws.Range("N122").Value = -67889.99800000001

$ws.Range("H125").Value = 84294.5
$ws.Range("J125").Value = 1034.1666
$ws.Range("L125").Value = 9307.499400000001
$ws.Range("N125").Value = -14227.4994

$ws.Range("H132").Value = 3647.9
$ws.Range("I132").Value = 3768.1177
$ws.Range("J132").Value = 2966.6667
$ws.Range("K132").Value = 11304.3531
$ws.Range("L132").Value = 8900.000100000001
$ws.Range("M132").Value = -8774.3531
$ws.Range("N132").Value = -13960.0001


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4440.8115
$ws.Range("I132").Value = 5741.393
$ws.Range("J132").Value = 2984.16
$ws.Range("K132").Value = 17224.179
$ws.Range("L132").Value = 8952.48
$ws.Range("M132").Value = -14694.179
$ws.Range("N132").Value = -14012.48


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 61753.117
$ws.Range("I20").Value = 84967
$ws.Range("J20").Value = 6039.8
$ws.Range("K20").Value = 84967
$ws.Range("L20").Value = 6039.8
$ws.Range("M20").Value = -84720
$ws.Range("N20").Value = -6533.8


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2254.1428
$ws.Range("I31").Value = 1328.6666
$ws.Range("J31").Value = 3920
$ws.Range("K31").Value = 1328.6666
$ws.Range("L31").Value = 3920
$ws.Range("M31").Value = -1033.6666
$ws.Range("N31").Value = -4510

$ws.Range("H34").Value = 2254.1428
$ws.Range("I34").Value = 1328.6666
$ws.Range("J34").Value = 3920
$ws.Range("K34").Value = 1328.6666
$ws.Range("L34").Value = 3920
$ws.Range("M34").Value = -1126.6666
$ws.Range("N34").Value = -4324

$ws.Range("H105").Value = 1711.625
$ws.Range("I105").Value = 1713.95
$ws.Range("J105").Value = 1700
$ws.Range("K105").Value = 1713.95
$ws.Range("L105").Value = 1700
$ws.Range("M105").Value = 33.04999999999995
$ws.Range("N105").Value = -5194


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 621.4286
$ws.Range("I47").Value = 516.6667
$ws.Range("J47").Value = 700
$ws.Range("K47").Value = 1550.0001
$ws.Range("L47").Value = 2100
$ws.Range("M47").Value = -1119.0001
$ws.Range("N47").Value = -2962

$ws.Range("H48").Value = 5282.875
$ws.Range("J48").Value = 5282.875
$ws.Range("L48").Value = 15848.625
$ws.Range("N48").Value = -16348.625

$ws.Range("H49").Value = 4000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 4000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 12000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -12312

$ws.Range("H81").Value = 6137.923
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 6566.0835
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 19698.2505
$ws.Range("M81").Value = -1877
$ws.Range("N81").Value = -21944.2505

$ws.Range("H84").Value = 6137.923
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 6566.0835
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 59094.7515
$ws.Range("M84").Value = -3384
$ws.Range("N84").Value = -70326.7515

$ws.Range("H113").Value = 697.3333
$ws.Range("I113").Value = 618.3182
$ws.Range("J113").Value = 855.36365
$ws.Range("K113").Value = 1854.9546
$ws.Range("L113").Value = 2566.09095
$ws.Range("M113").Value = 315.0454
$ws.Range("N113").Value = -6906.09095

$ws.Range("H134").Value = 5363.3125
$ws.Range("I134").Value = 3108
$ws.Range("J134").Value = 6388.4546
$ws.Range("K134").Value = 9324
$ws.Range("L134").Value = 19165.3638
$ws.Range("M134").Value = -4254
$ws.Range("N134").Value = -29305.3638


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 19000
$ws.Range("J39").Value = 19000
$ws.Range("L39").Value = 19000
$ws.Range("N39").Value = -20064

$ws.Range("H70").Value = 6097.4634
$ws.Range("I70").Value = 5720.433
$ws.Range("J70").Value = 7125.727
$ws.Range("K70").Value = 5720.433
$ws.Range("L70").Value = 7125.727
$ws.Range("M70").Value = -5450.433
$ws.Range("N70").Value = -7665.727

$ws.Range("H73").Value = 6097.4634
$ws.Range("I73").Value = 5720.433
$ws.Range("J73").Value = 7125.727
$ws.Range("K73").Value = 5720.433
$ws.Range("L73").Value = 7125.727
$ws.Range("M73").Value = -4784.433
$ws.Range("N73").Value = -8997.726999999999


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3532.2
$ws.Range("I40").Value = 4409.75
$ws.Range("J40").Value = 2529.2856
$ws.Range("K40").Value = 4409.75
$ws.Range("L40").Value = 2529.2856
$ws.Range("M40").Value = -4273.75
$ws.Range("N40").Value = -2801.2856

$ws.Range("H93").Value = 1864
$ws.Range("I93").Value = 1298
$ws.Range("J93").Value = 1958.3334
$ws.Range("K93").Value = 1298
$ws.Range("L93").Value = 1958.3334
$ws.Range("M93").Value = -50
$ws.Range("N93").Value = -4454.3334

$ws.Range("H132").Value = 4797.1113
$ws.Range("I132").Value = 4089.8667
$ws.Range("J132").Value = 8333.333000000001
$ws.Range("K132").Value = 12269.6001
$ws.Range("L132").Value = 24999.999
$ws.Range("M132").Value = -9739.6001
$ws.Range("N132").Value = -30059.999


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 53006
$ws.Range("I24").Value = 27500
$ws.Range("K24").Value = 27500
$ws.Range("M24").Value = -27270

$ws.Range("H121").Value = 30166.928
$ws.Range("J121").Value = 30166.928
$ws.Range("L121").Value = 30166.928
$ws.Range("N121").Value = -33660.928

$ws.Range("H136").Value = 2304.4583
$ws.Range("I136").Value = 2017.7368
$ws.Range("J136").Value = 3394
$ws.Range("K136").Value = 6053.2104
$ws.Range("L136").Value = 10182
$ws.Range("M136").Value = -3503.2104
$ws.Range("N136").Value = -15282

